$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 12:16"

# --- Madrid (row 4) gets updated totals ---
$ws.Range("B4").Value = 12352
$ws.Range("C4").Value = 2291
$ws.Range("D4").Value = 8526
$ws.Range("E4").Value = 1535

# --- Navarra / Araba-Alava swap order (rows 8-9) with fresh Navarra stats ---
$ws.Range("A8").Value = "Navarra"
$ws.Range("B8").Value = 1014
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 972
$ws.Range("E8").Value = 31

$ws.Range("A9").Value = "Araba/Alava"
$ws.Range("B9").Value = 1009
$ws.Range("C9").Value = 283
$ws.Range("D9").Value = 939
$ws.Range("E9").Value = 70

# --- La Rioja (row 10) gets updated totals ---
$ws.Range("B10").Value = 802
$ws.Range("C10").Value = 24
$ws.Range("D10").Value = 748
$ws.Range("E10").Value = 30

# --- Cantabria moves up to row 19 (before Salamanca) with fresh stats; ---
# --- Salamanca, Murcia, Gipuzkoa, Granada, Sevilla, Valladolid each shift down one row ---
$ws.Range("A19").Value = "Cantabria"
$ws.Range("B19").Value = 425
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 404
$ws.Range("E19").Value = 9

$ws.Range("A20").Value = "Salamanca"
$ws.Range("B20").Value = 404
$ws.Range("C20").Value = 22
$ws.Range("D20").Value = 358
$ws.Range("E20").Value = 24

$ws.Range("A21").Value = "Murcia"
$ws.Range("B21").Value = 385
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 381
$ws.Range("E21").Value = 3

$ws.Range("A22").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B22").Value = 380
$ws.Range("C22").Value = 283
$ws.Range("D22").Value = 365
$ws.Range("E22").Value = 15

$ws.Range("A23").Value = "Granada"
$ws.Range("B23").Value = 374
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 357
$ws.Range("E23").Value = 17

$ws.Range("A24").Value = "Sevilla"
$ws.Range("B24").Value = 351
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 345
$ws.Range("E24").Value = 5

$ws.Range("A25").Value = "Valladolid"
$ws.Range("B25").Value = 349
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 318
$ws.Range("E25").Value = 14

# Row 26 (Burgos) is unchanged: 336 / 29 / 289 / 18

Write-Output "Edit complete"
